# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / "heading" / "Outstanding" columns
# one place to the right (N->O, O->P, P->Q), then leave that sheet active
# with the selection on R7 (mirrors the manual Excel "Insert Column" edit
# captured in the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N; existing columns N..P shift to O..Q.
$ws.Columns("N").Insert() | Out-Null

# The newly inserted column inherits the left neighbour's (M) width in the
# source workbook (10.7109375); ColumnWidth is expressed in characters, so
# back it out to match as closely as this engine's rounding allows.
$ws.Columns("N").ColumnWidth = 9.8325

# Make "Repayment schedule" the active sheet/tab, with R7 selected - this
# matches the workbook-level activeTab move (Transactions -> Repayment
# schedule) and the sheetView/selection changes in the diff.
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
